# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns,
# and fix two swapped-row pairs (coin name/link/price/volume) so the
# ranking order matches the latest snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column D ("Price") updates whose new value is unambiguous text already
# (multi-dot thousand separators) - can be written directly.
# ---------------------------------------------------------------------------
$safeD = @{
    "D2"  = "67.080.38"
    "D3"  = "3.866.01"
    "D7"  = "3.861.94"
    "D15" = "4.473.28"
    "D17" = "3.859.30"
    "D20" = "67.138.53"
    "D39" = "0.0₃0763"
}
foreach ($addr in $safeD.Keys) {
    $ws.Range($addr).Value = $safeD[$addr]
}

# ---------------------------------------------------------------------------
# Column D ("Price") updates whose new value looks like a plain number
# (e.g. "426.61"). Force the cell to keep its original plain-text storage
# by applying a text number format before writing the value, otherwise
# Excel would silently convert it into a numeric value.
# ---------------------------------------------------------------------------
$textForceD = @{
    "D5"  = "426.61"
    "D6"  = "130.83"
    "D8"  = "0.612"
    "D13" = "40.98"
    "D14" = "10.28"
    "D19" = "19.96"
    "D21" = "1.09"
    "D22" = "413.24"
    "D23" = "14.96"
    "D24" = "84.25"
    "D25" = "3.03"
    "D26" = "37.79"
    "D27" = "10.00"
    "D28" = "3.25"
    "D30" = "8.90"
    "D31" = "725.83"
    "D35" = "0.999"
    "D36" = "39.06"
    "D37" = "0.152"
    "D38" = "56.00"
    "D41" = "0.0461"
    "D44" = "3.22"
    "D45" = "0.135"
    "D46" = "3.37"
    "D48" = "142.35"
    "D49" = "2.04"
    "D50" = "2.58"
    "D51" = "2.81"
}
foreach ($addr in $textForceD.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForceD[$addr]
}

# ---------------------------------------------------------------------------
# Column E ("Volume(1h)") updates - always padded with two leading and
# trailing spaces, so Excel keeps them as plain text automatically.
# ---------------------------------------------------------------------------
$eChanges = @{
    "E2"  = "  +7.55%  "
    "E3"  = "  +13.05%  "
    "E4"  = "  -0.16%  "
    "E5"  = "  +11.46%  "
    "E6"  = "  +8.71%  "
    "E7"  = "  +9.12%  "
    "E8"  = "  +6.61%  "
    "E9"  = "  -0.16%  "
    "E10" = "  +10.99%  "
    "E11" = "  +14.07%  "
    "E12" = "  +22.09%  "
    "E13" = "  +8.43%  "
    "E14" = "  +15.10%  "
    "E15" = "  +13.98%  "
    "E16" = "  +32.99%  "
    "E17" = "  +12.89%  "
    "E18" = "  +1.39%  "
    "E19" = "  +10.16%  "
    "E20" = "  +7.81%  "
    "E21" = "  +9.71%  "
    "E22" = "  +8.27%  "
    "E23" = "  +12.67%  "
    "E24" = "  +7.36%  "
    "E25" = "  +11.06%  "
    "E26" = "  +17.70%  "
    "E27" = "  +17.06%  "
    "E28" = "  +12.70%  "
    "E29" = "  +2.94%  "
    "E30" = "  +36.35%  "
    "E31" = "  +13.23%  "
    "E32" = "  +18.74%  "
    "E33" = "  +15.13%  "
    "E34" = "  +8.22%  "
    "E35" = "  -0.12%  "
    "E36" = "  +10.25%  "
    "E37" = "  +3.72%  "
    "E38" = "  +5.21%  "
    "E39" = "  +30.04%  "
    "E40" = "  +34.64%  "
    "E41" = "  +9.65%  "
    "E42" = "  +13.15%  "
    "E43" = "  +0.80%  "
    "E44" = "  +11.76%  "
    "E45" = "  +5.63%  "
    "E46" = "  +13.92%  "
    "E47" = "  +16.26%  "
    "E48" = "  +5.01%  "
    "E49" = "  +7.94%  "
    "E50" = "  +8.27%  "
    "E51" = "  +9.90%  "
}
foreach ($addr in $eChanges.Keys) {
    $ws.Range($addr).Value = $eChanges[$addr]
}

# ---------------------------------------------------------------------------
# Three ranking swaps where the relative order of two coins changed, so the
# Coin name (B) and Link (C) columns need to be swapped too:
#   row 14/15 - Polkadot <-> WrappedliquidstakedEther2.0
#   row 44/45 - ApeXProtocol <-> Stellar
#   row 50/51 - WEMIXToken <-> Stacks
# ---------------------------------------------------------------------------
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"

$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
